$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D6 currently holds text "4:36"; change it to a numeric value 436
$ws.Range("D6").Value = 436

# E6 currently holds text "17:05"; keep the same text value (it stays a string)
$ws.Range("E6").Value = "17:05"
